$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# --- Column width updates (widened to fit new hyperlink/file-name columns) ---
$overview.Columns.Item(5).ColumnWidth = 29.9777047293527
$overview.Columns.Item(6).ColumnWidth = 29.9777047293527

$zhcn.Columns.Item(3).ColumnWidth = 29.9777047293527
$zhcn.Columns.Item(9).ColumnWidth = 40
$zhcn.Columns.Item(10).ColumnWidth = 40

$dede.Columns.Item(3).ColumnWidth = 29.9777047293527
$dede.Columns.Item(9).ColumnWidth = 40
$dede.Columns.Item(10).ColumnWidth = 40

# --- Status column: "Ready for handoff" -> "Handed back: in sync with en-US" ---
$zhcn.Range("C2").Value = "Handed back: in sync with en-US"
$zhcn.Range("C3").Value = "Handed back: in sync with en-US"
$dede.Range("C2").Value = "Handed back: in sync with en-US"
$dede.Range("C3").Value = "Handed back: in sync with en-US"

# --- zh-cn: Latest Target File (I) + Latest Handback File (J) ---
$zhcn.Range("J2").Value = "8315e09e-2af9-4327-a12d-ac5760e73a7f.7315c86a1ee601c1b66d28aeae6ec5dabb543a56.zh-cn.xlf"
$zhcn.Range("J3").Value = "b9747646-22d8-4a72-b1f8-4868631950d6.b1a93591cbe06d96c1ac56d822548f6042eaa4e5.zh-cn.xlf"

$zhcn.Hyperlinks.Add($zhcn.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d21fcc7dbe8f7dbb81d8af3ac475e16f0d63a70f/e2e/8315e09e-2af9-4327-a12d-ac5760e73a7f.md", [Type]::Missing, [Type]::Missing, "8315e09e-2af9-4327-a12d-ac5760e73a7f.md") | Out-Null
$zhcn.Hyperlinks.Add($zhcn.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d21fcc7dbe8f7dbb81d8af3ac475e16f0d63a70f/e2e/b9747646-22d8-4a72-b1f8-4868631950d6.md", [Type]::Missing, [Type]::Missing, "b9747646-22d8-4a72-b1f8-4868631950d6.md") | Out-Null

# --- de-de: Latest Target File (I) + Latest Handback File (J) + Latest Handback DateTime (K) ---
$dede.Range("J2").Value = "8315e09e-2af9-4327-a12d-ac5760e73a7f.7315c86a1ee601c1b66d28aeae6ec5dabb543a56.de-de.xlf"
$dede.Range("J3").Value = "b9747646-22d8-4a72-b1f8-4868631950d6.b1a93591cbe06d96c1ac56d822548f6042eaa4e5.de-de.xlf"
$dede.Range("K2").Value = "2016-08-16 18:56:28"
$dede.Range("K3").Value = "2016-08-16 18:56:28"

$dede.Hyperlinks.Add($dede.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d21fcc7dbe8f7dbb81d8af3ac475e16f0d63a70f/e2e/8315e09e-2af9-4327-a12d-ac5760e73a7f.md", [Type]::Missing, [Type]::Missing, "8315e09e-2af9-4327-a12d-ac5760e73a7f.md") | Out-Null
$dede.Hyperlinks.Add($dede.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d21fcc7dbe8f7dbb81d8af3ac475e16f0d63a70f/e2e/b9747646-22d8-4a72-b1f8-4868631950d6.md", [Type]::Missing, [Type]::Missing, "b9747646-22d8-4a72-b1f8-4868631950d6.md") | Out-Null

# --- zh-cn: Latest Handback DateTime (K) text refresh (was placeholder epoch date) ---
$zhcn.Range("K2").Value = "2016-08-16 18:56:21"
$zhcn.Range("K3").Value = "2016-08-16 18:56:21"
